# Updated Power Pins and Header Layout
$wb = $excel.ActiveWorkbook

# --- "Links" sheet: move the selection (no longer the active/selected tab) ---
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Range("C15").Select() | Out-Null

# --- "Relevant Parts" sheet: new header/pin layout data ---
$wsParts = $wb.Worksheets.Item("Relevant Parts")

$wsParts.Range("A1").Value  = "6 pin female"
$wsParts.Range("A2").Value  = "8 pin female"
$wsParts.Range("A3").Value  = "12 pin female"
$wsParts.Range("A4").Value  = "6 pin male"
$wsParts.Range("A5").Value  = "8 pin male"
$wsParts.Range("A6").Value  = "12 pin male"
$wsParts.Range("A7").Value  = "10k resistor 0603"
$wsParts.Range("A8").Value  = "RGB LED"
$wsParts.Range("A9").Value  = "RGB Resistor"
$wsParts.Range("A10").Value = "PWR LED"
$wsParts.Range("A11").Value = "PWR Resistor"
$wsParts.Range("A12").Value = "Coin Cell Holder"
$wsParts.Range("A13").Value = "Coin Cell"
$wsParts.Range("A14").Value = "Switch"
$wsParts.Range("A15").Value = "Regulator"
$wsParts.Range("A16").Value = "4.7uF Cap_16V"
$wsParts.Range("A17").Value = "1uF Cap 0603"
$wsParts.Range("A18").Value = "0.1uF Cap 0603"
$wsParts.Range("A19").Value = "5v Regulator"
$wsParts.Range("A20").Value = "6 Pin IMU"

# Column B text values (order matters for shared-string table layout)
$wsParts.Range("B3").Value = "https://www.digikey.com/en/products/detail/amphenol-icc-fci/91601-412ALF/2713716"
$wsParts.Range("B2").Value = "https://www.digikey.com/en/products/detail/amphenol-icc-fci/91601-308LF/1535152"
$wsParts.Range("B1").Value = "https://www.digikey.com/en/products/detail/amphenol-icc-fci/91614-306ALF/1536169"
$wsParts.Range("B4").Value = "https://www.digikey.com/en/products/detail/amphenol-icc-fci/77311-401-06LF/1657857"
$wsParts.Range("B5").Value = "https://www.digikey.com/en/products/detail/amphenol-icc-fci/68002-208HLF/1878507"
$wsParts.Range("B6").Value = "(2 6 pins)"

# Hyperlinks (order matters for relationship-id allocation order)
$wsParts.Hyperlinks.Add($wsParts.Range("B2"), "https://www.digikey.com/en/products/detail/amphenol-icc-fci/91601-308LF/1535152") | Out-Null
$wsParts.Hyperlinks.Add($wsParts.Range("B4"), "https://www.digikey.com/en/products/detail/amphenol-icc-fci/77311-401-06LF/1657857") | Out-Null
$wsParts.Hyperlinks.Add($wsParts.Range("B3"), "https://www.digikey.com/en/products/detail/amphenol-icc-fci/91601-412ALF/2713716") | Out-Null
$wsParts.Hyperlinks.Add($wsParts.Range("B5"), "https://www.digikey.com/en/products/detail/amphenol-icc-fci/68002-208HLF/1878507") | Out-Null

# Auto-size column A to fit the new labels
$wsParts.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# This sheet becomes the active / selected tab, with this cell range selected
$wsParts.Activate() | Out-Null
$wsParts.Range("N13:N14").Select() | Out-Null
